# Generate Report for Handoff
# Replaces the two tracked files (af6f9aed... -> 99df6f75..., cfd1a895... -> ffffaa75c050...)
# and refreshes their handoff/handback status + timestamps across the
# Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldFile1 = "af6f9aed-1064-4e76-ac2d-1b234dce5b8b.md"
$oldFile2 = "cfd1a895-e6bd-4069-acff-fe2dd7aad653.md"
$newFile1 = "99df6f75-e167-486c-a91e-90f005fb0afe.md"
$newFile2 = "ffffaa75c050-105c-424f-8137-a52485764bb4.md"

$newZhXlf = "99df6f75-e167-486c-a91e-90f005fb0afe.091dd9ea77e38f935784c3c9bfeeb1546baf942e.zh-cn.xlf"
$newDeXlf = "99df6f75-e167-486c-a91e-90f005fb0afe.091dd9ea77e38f935784c3c9bfeeb1546baf942e.de-de.xlf"

$newStatus = "Ready for handoff"

# ============================== Overview sheet ==============================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = "e2e\" + $newFile1
$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = "e2e\" + $newFile2

$wsOverview.Range("E2:F3").Value = $newStatus
$wsOverview.Range("G2:G3").Value = "2016-08-18 11:05:50"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$B`$2") {
        $h.TextToDisplay = "e2e\" + $newFile1
    } elseif ($addr -eq "`$B`$3") {
        $h.TextToDisplay = "e2e\" + $newFile2
    }
}

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# ============================== zh-cn sheet ==============================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("C2:C3").Value = $newStatus
$wsZh.Range("G2:G3").Value = $newZhXlf
$wsZh.Range("H2:H3").Value = "2016-08-18 11:05:45"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("K2:K3").Value = "0001-01-01 00:00:00"

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = $newFile1
    } elseif ($addr -eq "`$A`$3") {
        $h.TextToDisplay = $newFile2
    } elseif ($addr -eq "`$I`$2" -or $addr -eq "`$I`$3") {
        $h.Delete()
    }
}

$wsZh.Range("I2:I3").Value = ""
$wsZh.Range("I2:I3").Style = "Normal"
$wsZh.Range("J2:J3").Value = ""

$wsZh.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ============================== de-de sheet ==============================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("C2:C3").Value = $newStatus
$wsDe.Range("G2:G3").Value = $newDeXlf
$wsDe.Range("H2:H3").Value = "2016-08-18 11:05:50"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("K2:K3").Value = "0001-01-01 00:00:00"

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = $newFile1
    } elseif ($addr -eq "`$A`$3") {
        $h.TextToDisplay = $newFile2
    } elseif ($addr -eq "`$I`$2" -or $addr -eq "`$I`$3") {
        $h.Delete()
    }
}

$wsDe.Range("I2:I3").Value = ""
$wsDe.Range("I2:I3").Style = "Normal"
$wsDe.Range("J2:J3").Value = ""

$wsDe.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
